$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) -- "想去人数" (interest counter) bumps only.
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$sheet1Updates = @{2=229; 3=55019; 4=1246; 5=17; 6=867; 7=28; 8=1139; 9=1405; 10=135; 11=45; 12=237; 13=408; 14=70; 15=42; 18=70; 19=5584; 21=5442; 22=9496; 23=119; 24=165; 25=160; 26=250; 27=464; 28=146; 29=115; 30=4303; 31=321}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# -------------------------------------------------------------------------
# Sheet 2: 演出 (Performance) -- the old row 2 ("东方乐典2024") disappeared
# from the listing, so every subsequent row shifts up by one; the running
# index in column A stays a plain 1..21 sequence; and the event that lands
# on row 8 ("EVA高桥洋子巡演") also had its counter tick up by one more.
# -------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows(2).EntireRow.Delete()
for ($i = 0; $i -lt 21; $i++) {
    $ws2.Cells.Item($i + 2, 1).Value = $i + 1
}
$ws2.Cells.Item(8, 6).Value = 1152

# -------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) -- two counter bumps, and the last row
# ("LoveLive! Series Asia Tour 2024") was removed entirely.
# -------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 787
$ws3.Cells.Item(3, 6).Value = 602
$ws3.Rows(5).EntireRow.Delete()

# -------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) -- mirrors the same counter bumps seen on
# sheets 1 and 3 (the rows it shares with them), no row restructuring.
# -------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$sheet4Updates = @{2=787; 3=602; 4=229; 5=1246; 6=17; 8=867; 9=28; 10=1139; 12=1405; 14=135; 15=237; 17=408; 18=70; 19=42; 23=70; 24=5584; 26=5442; 27=9496; 29=119; 30=165; 31=160; 32=250; 34=464; 38=146; 39=115; 40=4303; 47=321}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
